{"js": "// The document's first paragraph currently reads \"H?\" (split across two\n// runs \u2014 a leftover fragment of a heading question). The edit rewrites it\n// to the full question text, while keeping the first run's formatting\n// (Arial, bold, color 1F1F1F, sz 21, etc.) and collapsing the paragraph\n// down to a single run.\n\nconst body = context.document.body;\nconst firstParagraph = body.paragraphs.getFirst();\n\n// Replacing the whole paragraph range's text in one shot re-uses the\n// formatting of the paragraph's leading run and removes the extra run\n// that held the trailing \"?\", exactly matching the target OOXML.\nconst range = firstParagraph.getRange();\nrange.insertText(\n  \"Where are the best places for Students to live in Exeter, UK?\",\n  Word.InsertLocation.replace\n);\n\nawait context.sync();\n", "ps1": "# The document's first paragraph currently reads \"H?\" (split across two\n# runs -- a leftover fragment of a heading question). Rewrite it to the\n# full question text. Using Find/Replace scoped to that paragraph's range\n# merges the two runs into one, keeping the leading run's formatting\n# (Arial, bold, color 1F1F1F, sz 21, etc.) exactly like the target edit.\n\n$d = $word.ActiveDocument\n$p = $d.Paragraphs(1)\n$find = $p.Range.Find\n\n$find.ClearFormatting()\n$find.Execute(\n    \"H?\",                                                           # FindText\n    $false,                                                         # MatchCase\n    $false,                                                         # MatchWholeWord\n    $false,                                                         # MatchWildcards\n    $false,                                                         # MatchSoundsLike\n    $false,                                                         # MatchAllWordForms\n    $true,                                                          # Forward\n    1,                                                               # Wrap (wdFindContinue)\n    $false,                                                         # Format\n    \"Where are the best places for Students to live in Exeter, UK?\", # ReplaceWith\n    2                                                                # Replace (wdReplaceOne)\n)\n"}
